$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6746352193122216
$ws.Range("C2").Value = 0.9801841963585517

$ws.Range("B3").Value = 0.3676154600665907
$ws.Range("C3").Value = 0.9949194548047112

$ws.Range("B4").Value = 0.1534466317139522
$ws.Range("C4").Value = 0.9984184582858154

$ws.Range("B5").Value = 0.4516369355307997
$ws.Range("C5").Value = 0.9973097276373256
